$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value. Values must remain text
# (matching the source inlineStr cells), so we force the "@" text
# number format before assigning, then restore the default "Normal"
# style afterwards so no stray formatting is left on the cell.
$updates = @{
    "D2" = "265.46"
    "E2" = "1.58%"
    "D3" = "26.79"
    "E3" = "-1.19%"
    "D4" = "4.698"
    "E4" = "-0.14%"
    "D5" = "0.06083"
    "D6" = "6.733"
    "E6" = "0.65%"
    "D7" = "0.8507"
    "E7" = "0.04%"
    "D8" = "0.9085"
    "E8" = "-0.62%"
    "E9" = "0.05%"
    "D10" = "0.05044"
    "E10" = "8.64%"
    "D11" = "0.07097"
    "E11" = "0.06%"
    "D12" = "0.03150"
    "E12" = "-0.13%"
    "D13" = "0.09026"
    "E13" = "-0.20%"
    "D14" = "0.001529"
    "E14" = "-0.15%"
    "D15" = "0.0006093"
    "E15" = "-1.26%"
    "D16" = "0.005972"
    "E16" = "-1.09%"
    "E17" = "-0.03%"
    "D18" = "3.169"
    "E18" = "-0.22%"
    "E19" = "0.01%"
    "E21" = "-1.39%"
    "D22" = "4.086"
    "E22" = "-0.16%"
    "D23" = "0.04245"
    "E23" = "0.20%"
    "E24" = "-3.11%"
    "D25" = "0.004058"
    "E25" = "6.73%"
    "E26" = "0.07%"
    "E27" = "6.62%"
    "D40" = "0.03922"
    "E40" = "0.67%"
    "D41" = "0.1114"
    "E41" = "0.29%"
    "D42" = "0.004191"
    "E42" = "2.28%"
    "D43" = "0.002112"
    "E43" = "-3.31%"
    "D44" = "0.01160"
    "E44" = "-16.56%"
    "D45" = "0.00005111"
    "E45" = "-0.97%"
    "E46" = "0.07%"
    "D48" = "0.2580"
    "E48" = "53.86%"
    "E49" = "0.07%"
    "E50" = "0.07%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

